# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to remain plain text (matches the original inlineStr cells)
    # instead of letting Excel auto-convert numeric-looking strings to numbers.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "68.736.55"
Set-TextValue $ws.Range("E2") "  +1.42%  "
Set-TextValue $ws.Range("D3") "3.867.60"
Set-TextValue $ws.Range("E3") "  +0.82%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "602.33"
Set-TextValue $ws.Range("E5") "  +0.95%  "
Set-TextValue $ws.Range("D6") "171.77"
Set-TextValue $ws.Range("E6") "  +3.72%  "
Set-TextValue $ws.Range("D7") "3.868.34"
Set-TextValue $ws.Range("E7") "  +0.80%  "
Set-TextValue $ws.Range("E9") "  +1.10%  "
Set-TextValue $ws.Range("D10") "0.169"
Set-TextValue $ws.Range("E10") "  +2.99%  "
Set-TextValue $ws.Range("E11") "  +4.28%  "
Set-TextValue $ws.Range("E12") "  +1.76%  "
Set-TextValue $ws.Range("E13") "  +16.40%  "
Set-TextValue $ws.Range("D14") "37.28"
Set-TextValue $ws.Range("E14") "  +1.51%  "
Set-TextValue $ws.Range("D15") "4.515.22"
Set-TextValue $ws.Range("E15") "  +0.81%  "
Set-TextValue $ws.Range("D16") "3.877.95"
Set-TextValue $ws.Range("E16") "  +0.81%  "
Set-TextValue $ws.Range("D17") "68.781.83"
Set-TextValue $ws.Range("E17") "  +1.51%  "
Set-TextValue $ws.Range("D18") "18.40"
Set-TextValue $ws.Range("E18") "  +0.91%  "
Set-TextValue $ws.Range("D19") "7.41"
Set-TextValue $ws.Range("E19") "  +0.22%  "
Set-TextValue $ws.Range("E20") "  +0.70%  "
Set-TextValue $ws.Range("D21") "11.10"
Set-TextValue $ws.Range("E21") "  +4.16%  "
Set-TextValue $ws.Range("D22") "473.42"
Set-TextValue $ws.Range("E22") "  +1.59%  "
Set-TextValue $ws.Range("E23") "  +0.99%  "
Set-TextValue $ws.Range("E24") "  +2.59%  "
Set-TextValue $ws.Range("D25") "83.87"
Set-TextValue $ws.Range("E25") "  +1.09%  "
Set-TextValue $ws.Range("D26") "2.27"
Set-TextValue $ws.Range("E26") "  +3.73%  "
Set-TextValue $ws.Range("D27") "12.24"
Set-TextValue $ws.Range("E27") "  +1.21%  "
Set-TextValue $ws.Range("D28") "10.48"
Set-TextValue $ws.Range("E28") "  +5.55%  "
Set-TextValue $ws.Range("E29") "  +0.05%  "
Set-TextValue $ws.Range("E30") "  +1.11%  "
Set-TextValue $ws.Range("D31") "4.018.47"
Set-TextValue $ws.Range("E31") "  +0.88%  "
Set-TextValue $ws.Range("D32") "7.79"
Set-TextValue $ws.Range("E32") "  +1.74%  "
Set-TextValue $ws.Range("E33") "  +1.30%  "
Set-TextValue $ws.Range("D34") "31.37"
Set-TextValue $ws.Range("E34") "  +1.89%  "
Set-TextValue $ws.Range("D35") "9.43"
Set-TextValue $ws.Range("E35") "  +1.82%  "
Set-TextValue $ws.Range("D36") "3.832.25"
Set-TextValue $ws.Range("E36") "  +0.67%  "
Set-TextValue $ws.Range("D37") "3.98"
Set-TextValue $ws.Range("E37") "  +24.12%  "
Set-TextValue $ws.Range("D38") "0.106"
Set-TextValue $ws.Range("E38") "  +1.32%  "
Set-TextValue $ws.Range("E39") "  +1.72%  "
Set-TextValue $ws.Range("E40") "  +2.10%  "
Set-TextValue $ws.Range("E41") "  +0.64%  "
Set-TextValue $ws.Range("E42") "  +0.19%  "
Set-TextValue $ws.Range("D43") "0.321"
Set-TextValue $ws.Range("E43") "  +3.30%  "
Set-TextValue $ws.Range("E44") "  +15.78%  "
Set-TextValue $ws.Range("D45") "2.01"
Set-TextValue $ws.Range("E45") "  +1.95%  "
Set-TextValue $ws.Range("B46") "Cosmos"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D46") "8.82"
Set-TextValue $ws.Range("E46") "  +3.83%  "
Set-TextValue $ws.Range("B47") "USDe"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D47") "1.00"
Set-TextValue $ws.Range("E47") "  +0.03%  "
Set-TextValue $ws.Range("D48") "420.25"
Set-TextValue $ws.Range("E48") "  -0.65%  "
Set-TextValue $ws.Range("D49") "46.51"
Set-TextValue $ws.Range("E49") "  -1.23%  "
Set-TextValue $ws.Range("D50") "142.42"
Set-TextValue $ws.Range("E50") "  -0.48%  "
Set-TextValue $ws.Range("D51") "0.0360"
Set-TextValue $ws.Range("E51") "  +1.64%  "

Write-Host "Applied 86 cell updates"
